$wb = $excel.ActiveWorkbook

# The same correction needs to be applied on both the "All" sheet and the
# "Charge" sheet (they share the same underlying analyzed data for these rows).
$sheetNames = @("All", "Charge")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Row 3 updates
    $ws.Range("B3").Value = 0.02621867056820558
    $ws.Range("C3").Value = 0.9372963394055864
    $ws.Range("D3").Value = 0.03197733300429977
    $ws.Range("O3").Value = 112.4541184258475
    $ws.Range("P3").Value = ""

    # Row 4 updates
    $ws.Range("B4").Value = ""
    $ws.Range("C4").Value = 0.9146849503939857
    $ws.Range("D4").Value = ""
    $ws.Range("I4").Value = 302.7936205788868
    $ws.Range("J4").Value = ""
}
